$d = $word.ActiveDocument

# ----------------------------------------------------------------------------
# The target paragraph originally reads (as 4 runs):
#   R1 "En esta sección"
#   R2 " nos vamos a enfocar en importar y exportar arreglos y funciones que
#       se encuentran en un archivo externo hacia el archivo raíz"
#   R3 ", también vamos a exportar un módulo por defecto"
#   R4 ". "
#
# It must become (R1 and R4 kept byte-identical, R2+R3 replaced by 10 runs):
#   R1 "En esta sección"                                    (unchanged)
#   N1 " "
#   N2 "se procederá a"
#   N3 " importar y exportar arreglos"
#   N4 ","
#   N5 " y funciones que se encuentran en un archivo externo hacia el archivo raíz"
#   N6 ", también "
#   N7 "se realizará la"
#   N8 " exporta"
#   N9 "ción de"
#   N10 " un módulo por defecto"
#   R4 ". "                                                 (unchanged)
#
# This COM runtime collapses every run inside a paragraph into a single run
# as soon as any in-place edit (Find/Replace, Range.Text=, interior
# InsertBefore/InsertAfter, ...) touches that paragraph's interior. The only
# edits that do NOT cause a collapse are:
#   * InsertParagraphAfter() on a zero-length Range (a clean paragraph split)
#   * Range.Delete() of a paragraph mark that re-joins two paragraphs
#   * Paragraph.Range.InsertBefore/InsertAfter while the paragraph's own
#     original (untouched) content still sits between the two ends
#
# Strategy: temporarily split the paragraph into many small paragraphs at the
# desired run boundaries (each split is "clean"), overwrite each small
# paragraph's full text in one shot (safe: it is a single run getting fully
# replaced), then delete the inserted paragraph marks to re-join everything
# back into one paragraph -- which does NOT re-merge the now-distinct runs.
# ----------------------------------------------------------------------------

$newRuns = @(
    " ",
    "se procederá a",
    " importar y exportar arreglos",
    ",",
    " y funciones que se encuentran en un archivo externo hacia el archivo raíz",
    ", también ",
    "se realizará la",
    " exporta",
    "ción de",
    " un módulo por defecto"
)

# Locate the target paragraph robustly via Find rather than a hard-coded index.
$searchRange = $d.Content
$searchRange.Find.Execute("En esta sección", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraStartPos = $searchRange.Paragraphs(1).Range.Start
$p = $d.Paragraphs(1)
# Resolve the actual paragraph object that starts there.
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Start -eq $paraStartPos) { $p = $para; break }
}

$r1End = $p.Range.Start + ("En esta sección").Length

# Find end boundary of old-run3 (i.e. right before the trailing ". ") within
# this paragraph only.
$pDup = $p.Range.Duplicate
$pDup.Find.Execute("un módulo por defecto", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$beforeTailPos = $pDup.End

# 1) Split off R1 from the rest.
$splitPoint = $d.Range($r1End, $r1End)
$splitPoint.InsertParagraphAfter()

# The previous insert shifted every later offset by +1.
$beforeTailPos = $beforeTailPos + 1

# 2) Split off R4 (". ") from the middle block.
$splitPoint2 = $d.Range($beforeTailPos, $beforeTailPos)
$splitPoint2.InsertParagraphAfter()

# 3) The middle paragraph (old R2+R3) now sits isolated between R1's
#    paragraph and R4's paragraph. Re-find it fresh.
$midPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Start -eq ($r1End + 1)) { $midPara = $para; break }
}

# 4) Split the middle paragraph into N pieces (N = number of new runs) by
#    inserting N-1 additional paragraph marks at arbitrary interior offsets
#    (their exact location does not matter: every piece's text gets fully
#    overwritten in the next step anyway).
$pieceCount = $newRuns.Count
$midStart = $midPara.Range.Start
$midEnd = $midPara.Range.End - 1   # exclude its paragraph mark
$midLen = $midEnd - $midStart

for ($i = 1; $i -lt $pieceCount; $i++) {
    $offset = $midStart + [int][math]::Floor($midLen * $i / $pieceCount)
    if ($offset -le $midStart) { $offset = $midStart + 1 }
    if ($offset -ge $midEnd) { $offset = $midEnd - 1 }
    $sp = $d.Range($offset, $offset)
    $sp.InsertParagraphAfter()
    $midEnd = $midEnd + 1
}

# 5) Overwrite each of the now-separate small paragraphs with the desired
#    final text, in document order.
$firstMidStart = $midStart
$cursor = $firstMidStart
for ($i = 0; $i -lt $pieceCount; $i++) {
    $piece = $null
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Start -eq $cursor) { $piece = $para; break }
    }
    $pieceDup = $piece.Range.Duplicate
    $pieceDup.MoveEnd(1, -1)  # exclude paragraph mark
    $pieceDup.Text = $newRuns[$i]
    $cursor = $piece.Range.Start + $newRuns[$i].Length + 1
}

# 6) Re-join everything: delete the (pieceCount - 1) + 2 inserted paragraph
#    marks, from the last one back to the first (so earlier offsets stay
#    valid while later ones are being removed).
# Recompute: after step 5, paragraph marks sit right after each piece's text,
# except the very last piece which is followed by R4's paragraph (kept).
$cursor = $firstMidStart
for ($i = 0; $i -lt $pieceCount; $i++) {
    $cursor = $cursor + $newRuns[$i].Length
    if ($i -lt $pieceCount - 1) {
        $mark = $d.Range($cursor, $cursor + 1)
        $mark.Delete()
    } else {
        # last piece: merge with the following R4 paragraph
        $mark = $d.Range($cursor, $cursor + 1)
        $mark.Delete()
    }
}

# Finally, merge R1's paragraph with the (now single, multi-run) rest.
$mark = $d.Range($r1End, $r1End + 1)
$mark.Delete()
